# Commit: "Added timeout to redirect from create new listing/listing edit"
#
# Changes applied:
#  1. D13: append a new sentence to the "Bugfixes, clean up/formatting" log entry
#     describing the new redirect timeout.
#  2. B13: hours logged for that entry bumped from 5 to 8 (B30 total recalculates
#     automatically via its SUM formula).
#  3. Row 13 height grows (autofit for the extra wrapped line) from 51 to 75.75.
#  4. The sheet's active-cell selection moves from D14 to D20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1 & 2: update the log entry text and hours for row 13
$ws.Range("D13").Value = "Removed the ability for users to edit a page when a specific url is typed into the browser.`nFixed an issue where page redirection broke when a number without decimal speration was given.`nAdded a timeout to page redirect to allow image to be uploaded to API before user enters details page.`nCleaned up my code."
$ws.Range("B13").Value = 8

# 3: grow row 13 to fit the extra wrapped line of text
$ws.Rows.Item(13).RowHeight = 75.75

# 4: move the selected/active cell to D20
$ws.Range("D20").Select()
